# Correção nos dados e início da análise PNAD 2009
# Fix header labels: the "unnamed: X_level_1" placeholders produced by the
# pandas MultiIndex export should read "total" instead, matching the
# sibling "total" column header already present in C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
